# Refresh the crypto price/volume table in place, cell by cell, to match the
# latest scrape (GitHub Actions run). Also fixes the Aptos/Elrond row order
# (rows 48-49 swap places).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference, new text, and whether it must be forced to
# Text format first (otherwise Excel would auto-convert a numeric-looking
# "Price" string like "0.9998" into a float and mangle it).
$updates = @(
    @{Cell='D2'; Value='30.393.64'; AsText=$false}
    @{Cell='E2'; Value='  -0.74%  '; AsText=$false}
    @{Cell='D3'; Value='1.871.94'; AsText=$false}
    @{Cell='D4'; Value='0.9998'; AsText=$true}
    @{Cell='D5'; Value='244.12'; AsText=$true}
    @{Cell='E5'; Value='  -1.77%  '; AsText=$false}
    @{Cell='D6'; Value='0.9998'; AsText=$true}
    @{Cell='E6'; Value='  -0.04%  '; AsText=$false}
    @{Cell='D7'; Value='0.4713'; AsText=$true}
    @{Cell='E7'; Value='  -1.10%  '; AsText=$false}
    @{Cell='D8'; Value='0.2888'; AsText=$true}
    @{Cell='E8'; Value='  -1.43%  '; AsText=$false}
    @{Cell='D9'; Value='0.06506'; AsText=$true}
    @{Cell='E9'; Value='  -0.42%  '; AsText=$false}
    @{Cell='D10'; Value='22.00'; AsText=$true}
    @{Cell='E10'; Value='  +0.01%  '; AsText=$false}
    @{Cell='D11'; Value='0.07790'; AsText=$true}
    @{Cell='E11'; Value='  +0.70%  '; AsText=$false}
    @{Cell='D12'; Value='96.36'; AsText=$true}
    @{Cell='E12'; Value='  -0.41%  '; AsText=$false}
    @{Cell='D13'; Value='1.870.64'; AsText=$false}
    @{Cell='E13'; Value='  -0.42%  '; AsText=$false}
    @{Cell='D14'; Value='0.7253'; AsText=$true}
    @{Cell='E14'; Value='  -2.14%  '; AsText=$false}
    @{Cell='D15'; Value='5.151'; AsText=$true}
    @{Cell='E15'; Value='  -1.22%  '; AsText=$false}
    @{Cell='D16'; Value='283.04'; AsText=$true}
    @{Cell='E16'; Value='  +3.31%  '; AsText=$false}
    @{Cell='D17'; Value='30.382.86'; AsText=$false}
    @{Cell='E17'; Value='  -1.08%  '; AsText=$false}
    @{Cell='D18'; Value='13.04'; AsText=$true}
    @{Cell='E18'; Value='  -1.68%  '; AsText=$false}
    @{Cell='D19'; Value='0.000007526'; AsText=$true}
    @{Cell='E19'; Value='  -0.10%  '; AsText=$false}
    @{Cell='D21'; Value='2.112.58'; AsText=$false}
    @{Cell='D22'; Value='0.9995'; AsText=$true}
    @{Cell='E22'; Value='  -0.08%  '; AsText=$false}
    @{Cell='D23'; Value='5.268'; AsText=$true}
    @{Cell='E23'; Value='  +0.07%  '; AsText=$false}
    @{Cell='D24'; Value='6.247'; AsText=$true}
    @{Cell='E24'; Value='  +0.63%  '; AsText=$false}
    @{Cell='D25'; Value='163.86'; AsText=$true}
    @{Cell='E25'; Value='  -0.96%  '; AsText=$false}
    @{Cell='D26'; Value='9.077'; AsText=$true}
    @{Cell='E26'; Value='  -1.36%  '; AsText=$false}
    @{Cell='D27'; Value='18.81'; AsText=$true}
    @{Cell='E27'; Value='  -0.49%  '; AsText=$false}
    @{Cell='D28'; Value='1.885'; AsText=$true}
    @{Cell='E28'; Value='  -1.61%  '; AsText=$false}
    @{Cell='E29'; Value='  -1.37%  '; AsText=$false}
    @{Cell='D30'; Value='0.09624'; AsText=$true}
    @{Cell='E30'; Value='  -2.30%  '; AsText=$false}
    @{Cell='D31'; Value='1.488'; AsText=$true}
    @{Cell='E31'; Value='  -1.04%  '; AsText=$false}
    @{Cell='D32'; Value='4.240'; AsText=$true}
    @{Cell='E32'; Value='  -1.23%  '; AsText=$false}
    @{Cell='D33'; Value='4.131'; AsText=$true}
    @{Cell='E33'; Value='  +0.29%  '; AsText=$false}
    @{Cell='D34'; Value='0.04842'; AsText=$true}
    @{Cell='E34'; Value='  +0.16%  '; AsText=$false}
    @{Cell='E35'; Value='  -0.34%  '; AsText=$false}
    @{Cell='D36'; Value='0.6924'; AsText=$true}
    @{Cell='E36'; Value='  -0.55%  '; AsText=$false}
    @{Cell='E37'; Value='  -0.16%  '; AsText=$false}
    @{Cell='D38'; Value='0.01896'; AsText=$true}
    @{Cell='E38'; Value='  +0.91%  '; AsText=$false}
    @{Cell='D39'; Value='2.817'; AsText=$true}
    @{Cell='E39'; Value='  +1.96%  '; AsText=$false}
    @{Cell='D40'; Value='76.12'; AsText=$true}
    @{Cell='E40'; Value='  +3.64%  '; AsText=$false}
    @{Cell='D41'; Value='6.257'; AsText=$true}
    @{Cell='E41'; Value='  -0.23%  '; AsText=$false}
    @{Cell='D42'; Value='0.4238'; AsText=$true}
    @{Cell='E42'; Value='  -0.17%  '; AsText=$false}
    @{Cell='D43'; Value='1.937'; AsText=$true}
    @{Cell='E43'; Value='  -2.82%  '; AsText=$false}
    @{Cell='D44'; Value='0.9989'; AsText=$true}
    @{Cell='E44'; Value='  -0.16%  '; AsText=$false}
    @{Cell='D45'; Value='0.8291'; AsText=$true}
    @{Cell='E45'; Value='  -1.14%  '; AsText=$false}
    @{Cell='D46'; Value='101.01'; AsText=$true}
    @{Cell='E46'; Value='  -1.24%  '; AsText=$false}
    @{Cell='D47'; Value='9.673'; AsText=$true}
    @{Cell='E47'; Value='  +3.25%  '; AsText=$false}
    @{Cell='B48'; Value='Aptos'; AsText=$false}
    @{Cell='C48'; Value='https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; AsText=$false}
    @{Cell='D48'; Value='6.987'; AsText=$true}
    @{Cell='E48'; Value='  -0.64%  '; AsText=$false}
    @{Cell='B49'; Value='Elrond'; AsText=$false}
    @{Cell='C49'; Value='https://coinranking.com/coin/omwkOTglq+elrond-egld'; AsText=$false}
    @{Cell='D49'; Value='35.34'; AsText=$true}
    @{Cell='E49'; Value='  -0.35%  '; AsText=$false}
    @{Cell='D50'; Value='908.44'; AsText=$true}
    @{Cell='E50'; Value='  -0.22%  '; AsText=$false}
    @{Cell='D51'; Value='0.05729'; AsText=$true}
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    if ($u.AsText) {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $u.Value
}
